{"js": "const pairs = [\n  [\"54\u00d731=1674\", \"83\u00d730=2490\"],\n  [\"12\u00d735=420\", \"76\u00d760=4560\"],\n  [\"97\u00d726=2522\", \"44\u00d731=1364\"],\n  [\"74\u00d771=5254\", \"35\u00d747=1645\"],\n  [\"73\u00d779=5767\", \"98\u00d716=1568\"],\n  [\"63\u00d751=3213\", \"97\u00d781=7857\"],\n  [\"10\u00d736=360\", \"72\u00d795=6840\"],\n  [\"50\u00d799=4950\", \"24\u00d746=1104\"],\n  [\"65\u00d759=3835\", \"39\u00d710=390\"],\n  [\"29\u00d794=2726\", \"16\u00d757=912\"],\n  [\"18\u00d751=918\", \"84\u00d760=5040\"],\n  [\"35\u00d798=3430\", \"66\u00d714=924\"],\n  [\"13\u00d770=910\", \"79\u00d722=1738\"],\n  [\"36\u00d788=3168\", \"75\u00d756=4200\"],\n  [\"34\u00d795=3230\", \"99\u00d710=990\"],\n  [\"47\u00d770=3290\", \"46\u00d746=2116\"],\n  [\"90\u00d792=8280\", \"68\u00d765=4420\"],\n  [\"47\u00d754=2538\", \"93\u00d756=5208\"],\n  [\"25\u00d789=2225\", \"22\u00d748=1056\"],\n  [\"38\u00d775=2850\", \"91\u00d733=3003\"],\n  [\"96\u00d716=1536\", \"64\u00d783=5312\"],\n  [\"88\u00d777=6776\", \"51\u00d733=1683\"],\n  [\"28\u00d711=308\", \"34\u00d752=1768\"],\n  [\"47\u00d766=3102\", \"40\u00d766=2640\"],\n  [\"61\u00d798=5978\", \"90\u00d720=1800\"],\n  [\"100\u00d750=5000\", \"87\u00d775=6525\"],\n  [\"91\u00d741=3731\", \"33\u00d754=1782\"],\n  [\"87\u00d771=6177\", \"41\u00d739=1599\"],\n  [\"59\u00d779=4661\", \"81\u00d7100=8100\"],\n  [\"38\u00d763=2394\", \"90\u00d717=1530\"],\n  [\"86\u00d726=2236\", \"41\u00d760=2460\"],\n  [\"52\u00d769=3588\", \"92\u00d771=6532\"],\n  [\"57\u00d788=5016\", \"11\u00d780=880\"],\n  [\"21\u00d716=336\", \"87\u00d774=6438\"],\n  [\"22\u00d757=1254\", \"62\u00d729=1798\"],\n  [\"97\u00d717=1649\", \"27\u00d761=1647\"],\n  [\"26\u00d775=1950\", \"43\u00d739=1677\"],\n  [\"32\u00d720=640\", \"96\u00d744=4224\"],\n  [\"27\u00d723=621\", \"15\u00d794=1410\"],\n  [\"40\u00d717=680\", \"54\u00d794=5076\"],\n  [\"36\u00d727=972\", \"55\u00d752=2860\"],\n  [\"35\u00d741=1435\", \"41\u00d774=3034\"],\n  [\"85\u00d781=6885\", \"79\u00d735=2765\"],\n  [\"85\u00d727=2295\", \"20\u00d710=200\"],\n  [\"42\u00d746=1932\", \"13\u00d723=299\"],\n  [\"65\u00d745=2925\", \"89\u00d763=5607\"],\n  [\"53\u00d718=954\", \"16\u00d721=336\"],\n  [\"45\u00d762=2790\", \"82\u00d784=6888\"],\n  [\"26\u00d785=2210\", \"35\u00d777=2695\"],\n  [\"19\u00d716=304\", \"16\u00d714=224\"],\n  [\"22\u00d722=484\", \"44\u00d734=1496\"],\n  [\"82\u00d710=820\", \"56\u00d728=1568\"],\n  [\"39\u00d799=3861\", \"76\u00d790=6840\"],\n  [\"12\u00d761=732\", \"16\u00d792=1472\"],\n  [\"27\u00d799=2673\", \"60\u00d799=5940\"],\n  [\"41\u00d723=943\", \"56\u00d786=4816\"],\n  [\"57\u00d752=2964\", \"92\u00d772=6624\"],\n  [\"21\u00d771=1491\", \"80\u00d754=4320\"],\n  [\"17\u00d758=986\", \"44\u00d722=968\"],\n  [\"95\u00d742=3990\", \"88\u00d729=2552\"],\n  [\"84\u00d768=5712\", \"36\u00d783=2988\"],\n  [\"89\u00d736=3204\", \"48\u00d754=2592\"],\n  [\"16\u00d750=800\", \"19\u00d781=1539\"],\n  [\"97\u00d760=5820\", \"59\u00d769=4071\"],\n  [\"96\u00d736=3456\", \"98\u00d790=8820\"],\n  [\"90\u00d772=6480\", \"92\u00d721=1932\"],\n  [\"96\u00d713=1248\", \"88\u00d741=3608\"],\n  [\"18\u00d759=1062\", \"13\u00d723=299\"],\n  [\"98\u00d782=8036\", \"13\u00d724=312\"],\n  [\"13\u00d787=1131\", \"81\u00d773=5913\"],\n  [\"64\u00d765=4160\", \"77\u00d761=4697\"],\n  [\"74\u00d751=3774\", \"24\u00d728=672\"],\n  [\"14\u00d774=1036\", \"25\u00d712=300\"],\n  [\"85\u00d710=850\", \"28\u00d769=1932\"],\n  [\"74\u00d752=3848\", \"35\u00d782=2870\"],\n  [\"48\u00d711=528\", \"17\u00d719=323\"],\n  [\"82\u00d797=7954\", \"54\u00d7100=5400\"],\n  [\"67\u00d710=670\", \"64\u00d756=3584\"],\n  [\"51\u00d794=4794\", \"74\u00d774=5476\"],\n  [\"69\u00d749=3381\", \"10\u00d770=700\"],\n  [\"66\u00d739=2574\", \"55\u00d719=1045\"],\n  [\"91\u00d757=5187\", \"61\u00d782=5002\"],\n  [\"61\u00d767=4087\", \"11\u00d771=781\"],\n  [\"21\u00d743=903\", \"13\u00d739=507\"],\n  [\"78\u00d755=4290\", \"26\u00d766=1716\"],\n  [\"69\u00d797=6693\", \"27\u00d722=594\"],\n  [\"29\u00d787=2523\", \"83\u00d734=2822\"],\n  [\"39\u00d718=702\", \"11\u00d771=781\"],\n  [\"25\u00d778=1950\", \"83\u00d786=7138\"],\n  [\"70\u00d745=3150\", \"84\u00d774=6216\"],\n  [\"87\u00d735=3045\", \"83\u00d738=3154\"],\n  [\"13\u00d737=481\", \"90\u00d742=3780\"],\n  [\"26\u00d799=2574\", \"12\u00d716=192\"],\n  [\"81\u00d760=4860\", \"36\u00d777=2772\"],\n  [\"18\u00d724=432\", \"45\u00d792=4140\"],\n  [\"67\u00d772=4824\", \"38\u00d768=2584\"],\n  [\"17\u00d775=1275\", \"53\u00d725=1325\"],\n  [\"39\u00d796=3744\", \"74\u00d765=4810\"],\n  [\"18\u00d739=702\", \"18\u00d770=1260\"],\n  [\"37\u00d727=999\", \"40\u00d737=1480\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  ,@(\"54\u00d731=1674\", \"83\u00d730=2490\")\n  ,@(\"12\u00d735=420\", \"76\u00d760=4560\")\n  ,@(\"97\u00d726=2522\", \"44\u00d731=1364\")\n  ,@(\"74\u00d771=5254\", \"35\u00d747=1645\")\n  ,@(\"73\u00d779=5767\", \"98\u00d716=1568\")\n  ,@(\"63\u00d751=3213\", \"97\u00d781=7857\")\n  ,@(\"10\u00d736=360\", \"72\u00d795=6840\")\n  ,@(\"50\u00d799=4950\", \"24\u00d746=1104\")\n  ,@(\"65\u00d759=3835\", \"39\u00d710=390\")\n  ,@(\"29\u00d794=2726\", \"16\u00d757=912\")\n  ,@(\"18\u00d751=918\", \"84\u00d760=5040\")\n  ,@(\"35\u00d798=3430\", \"66\u00d714=924\")\n  ,@(\"13\u00d770=910\", \"79\u00d722=1738\")\n  ,@(\"36\u00d788=3168\", \"75\u00d756=4200\")\n  ,@(\"34\u00d795=3230\", \"99\u00d710=990\")\n  ,@(\"47\u00d770=3290\", \"46\u00d746=2116\")\n  ,@(\"90\u00d792=8280\", \"68\u00d765=4420\")\n  ,@(\"47\u00d754=2538\", \"93\u00d756=5208\")\n  ,@(\"25\u00d789=2225\", \"22\u00d748=1056\")\n  ,@(\"38\u00d775=2850\", \"91\u00d733=3003\")\n  ,@(\"96\u00d716=1536\", \"64\u00d783=5312\")\n  ,@(\"88\u00d777=6776\", \"51\u00d733=1683\")\n  ,@(\"28\u00d711=308\", \"34\u00d752=1768\")\n  ,@(\"47\u00d766=3102\", \"40\u00d766=2640\")\n  ,@(\"61\u00d798=5978\", \"90\u00d720=1800\")\n  ,@(\"100\u00d750=5000\", \"87\u00d775=6525\")\n  ,@(\"91\u00d741=3731\", \"33\u00d754=1782\")\n  ,@(\"87\u00d771=6177\", \"41\u00d739=1599\")\n  ,@(\"59\u00d779=4661\", \"81\u00d7100=8100\")\n  ,@(\"38\u00d763=2394\", \"90\u00d717=1530\")\n  ,@(\"86\u00d726=2236\", \"41\u00d760=2460\")\n  ,@(\"52\u00d769=3588\", \"92\u00d771=6532\")\n  ,@(\"57\u00d788=5016\", \"11\u00d780=880\")\n  ,@(\"21\u00d716=336\", \"87\u00d774=6438\")\n  ,@(\"22\u00d757=1254\", \"62\u00d729=1798\")\n  ,@(\"97\u00d717=1649\", \"27\u00d761=1647\")\n  ,@(\"26\u00d775=1950\", \"43\u00d739=1677\")\n  ,@(\"32\u00d720=640\", \"96\u00d744=4224\")\n  ,@(\"27\u00d723=621\", \"15\u00d794=1410\")\n  ,@(\"40\u00d717=680\", \"54\u00d794=5076\")\n  ,@(\"36\u00d727=972\", \"55\u00d752=2860\")\n  ,@(\"35\u00d741=1435\", \"41\u00d774=3034\")\n  ,@(\"85\u00d781=6885\", \"79\u00d735=2765\")\n  ,@(\"85\u00d727=2295\", \"20\u00d710=200\")\n  ,@(\"42\u00d746=1932\", \"13\u00d723=299\")\n  ,@(\"65\u00d745=2925\", \"89\u00d763=5607\")\n  ,@(\"53\u00d718=954\", \"16\u00d721=336\")\n  ,@(\"45\u00d762=2790\", \"82\u00d784=6888\")\n  ,@(\"26\u00d785=2210\", \"35\u00d777=2695\")\n  ,@(\"19\u00d716=304\", \"16\u00d714=224\")\n  ,@(\"22\u00d722=484\", \"44\u00d734=1496\")\n  ,@(\"82\u00d710=820\", \"56\u00d728=1568\")\n  ,@(\"39\u00d799=3861\", \"76\u00d790=6840\")\n  ,@(\"12\u00d761=732\", \"16\u00d792=1472\")\n  ,@(\"27\u00d799=2673\", \"60\u00d799=5940\")\n  ,@(\"41\u00d723=943\", \"56\u00d786=4816\")\n  ,@(\"57\u00d752=2964\", \"92\u00d772=6624\")\n  ,@(\"21\u00d771=1491\", \"80\u00d754=4320\")\n  ,@(\"17\u00d758=986\", \"44\u00d722=968\")\n  ,@(\"95\u00d742=3990\", \"88\u00d729=2552\")\n  ,@(\"84\u00d768=5712\", \"36\u00d783=2988\")\n  ,@(\"89\u00d736=3204\", \"48\u00d754=2592\")\n  ,@(\"16\u00d750=800\", \"19\u00d781=1539\")\n  ,@(\"97\u00d760=5820\", \"59\u00d769=4071\")\n  ,@(\"96\u00d736=3456\", \"98\u00d790=8820\")\n  ,@(\"90\u00d772=6480\", \"92\u00d721=1932\")\n  ,@(\"96\u00d713=1248\", \"88\u00d741=3608\")\n  ,@(\"18\u00d759=1062\", \"13\u00d723=299\")\n  ,@(\"98\u00d782=8036\", \"13\u00d724=312\")\n  ,@(\"13\u00d787=1131\", \"81\u00d773=5913\")\n  ,@(\"64\u00d765=4160\", \"77\u00d761=4697\")\n  ,@(\"74\u00d751=3774\", \"24\u00d728=672\")\n  ,@(\"14\u00d774=1036\", \"25\u00d712=300\")\n  ,@(\"85\u00d710=850\", \"28\u00d769=1932\")\n  ,@(\"74\u00d752=3848\", \"35\u00d782=2870\")\n  ,@(\"48\u00d711=528\", \"17\u00d719=323\")\n  ,@(\"82\u00d797=7954\", \"54\u00d7100=5400\")\n  ,@(\"67\u00d710=670\", \"64\u00d756=3584\")\n  ,@(\"51\u00d794=4794\", \"74\u00d774=5476\")\n  ,@(\"69\u00d749=3381\", \"10\u00d770=700\")\n  ,@(\"66\u00d739=2574\", \"55\u00d719=1045\")\n  ,@(\"91\u00d757=5187\", \"61\u00d782=5002\")\n  ,@(\"61\u00d767=4087\", \"11\u00d771=781\")\n  ,@(\"21\u00d743=903\", \"13\u00d739=507\")\n  ,@(\"78\u00d755=4290\", \"26\u00d766=1716\")\n  ,@(\"69\u00d797=6693\", \"27\u00d722=594\")\n  ,@(\"29\u00d787=2523\", \"83\u00d734=2822\")\n  ,@(\"39\u00d718=702\", \"11\u00d771=781\")\n  ,@(\"25\u00d778=1950\", \"83\u00d786=7138\")\n  ,@(\"70\u00d745=3150\", \"84\u00d774=6216\")\n  ,@(\"87\u00d735=3045\", \"83\u00d738=3154\")\n  ,@(\"13\u00d737=481\", \"90\u00d742=3780\")\n  ,@(\"26\u00d799=2574\", \"12\u00d716=192\")\n  ,@(\"81\u00d760=4860\", \"36\u00d777=2772\")\n  ,@(\"18\u00d724=432\", \"45\u00d792=4140\")\n  ,@(\"67\u00d772=4824\", \"38\u00d768=2584\")\n  ,@(\"17\u00d775=1275\", \"53\u00d725=1325\")\n  ,@(\"39\u00d796=3744\", \"74\u00d765=4810\")\n  ,@(\"18\u00d739=702\", \"18\u00d770=1260\")\n  ,@(\"37\u00d727=999\", \"40\u00d737=1480\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        Write-Output \"MISSING: $oldText\"\n    }\n}\n"}
